$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2 through 143)
# from serial date 45177 to 45178 (i.e. one day later).
$ws.Range("C2:C143").Value = 45178
